{"js": "// Smart/curly quote characters built from code points to stay encoding-safe.\nconst ldq = String.fromCharCode(0x201C); // \u201c\nconst rdq = String.fromCharCode(0x201D); // \u201d\n\nasync function replaceText(findText, newText) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 1. \"Go change its name then open the website\" -> \"Go change the name of the style.css file then open the website\"\nawait replaceText(\n  \"Go change its name then open the website\",\n  \"Go change the name of the style.css file then open the website\"\n);\n\n// 2. \"...uses the style.css file again.\" -> \"...uses the style.css file correctly again.\"\nawait replaceText(\n  \"uses the style.css file again.\",\n  \"uses the style.css file correctly again.\"\n);\n\n// 3. \"choose all files, \" -> 'choose \"All files\" instead of \"Text documents (*.txt)\", '\nawait replaceText(\n  \"choose all files, \",\n  \"choose \" + ldq + \"All files\" + rdq + \" instead of \" + ldq + \"Text documents (*.txt)\" + rdq + \", \"\n);\n", "ps1": "$d = $word.ActiveDocument\n\n# Smart/curly quote characters used below (avoid embedding literal non-ASCII\n# in the script source to stay encoding-safe across hosts).\n$ldq = [char]0x201C   # \"\n$rdq = [char]0x201D   # \"\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$replaceText, 2) | Out-Null\n}\n\n# 1. \"Go change its name then open the website\" -> \"Go change the name of the style.css file then open the website\"\nReplace-Text \"Go change its name then open the website\" \"Go change the name of the style.css file then open the website\"\n\n# 2. \"...uses the style.css file again.\" -> \"...uses the style.css file correctly again.\"\nReplace-Text \"uses the style.css file again.\" \"uses the style.css file correctly again.\"\n\n# 3. \"choose all files, \" -> 'choose \"All files\" instead of \"Text documents (*.txt)\", '\n$repl3 = \"choose \" + $ldq + \"All files\" + $rdq + \" instead of \" + $ldq + \"Text documents (*.txt)\" + $rdq + \", \"\nReplace-Text \"choose all files, \" $repl3\n"}
